$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.896799999999993
$ws.Range("A3").Value = -21.30570000000003
$ws.Range("D5").Value = -8.282199999999991
$ws.Range("A14").Value = -20.48659999999998
$ws.Range("A21").Value = -21.3022
$ws.Range("A23").Value = -21.29470000000002
$ws.Range("A25").Value = -22.41490000000003
